{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell with its new value.\n// Every original text value in this document is unique, so a plain\n// matchCase exact-text search/replace is safe and unambiguous.\nconst replacements = [\n  [\"2024-07-08 Monday\", \"2024-07-09 Tuesday\"],\n  [\"909\u00d77=6363\", \"614\u00d77=4298\"],\n  [\"887\u00d74=3548\", \"721\u00d79=6489\"],\n  [\"128\u00d77=896\", \"731\u00d72=1462\"],\n  [\"617\u00d78=4936\", \"384\u00d77=2688\"],\n  [\"559\u00d74=2236\", \"332\u00d73=996\"],\n  [\"150\u00d72=300\", \"790\u00d78=6320\"],\n  [\"584\u00d77=4088\", \"280\u00d74=1120\"],\n  [\"995\u00d73=2985\", \"384\u00d72=768\"],\n  [\"646\u00d78=5168\", \"370\u00d73=1110\"],\n  [\"490\u00d73=1470\", \"761\u00d78=6088\"],\n  [\"938\u00d72=1876\", \"743\u00d78=5944\"],\n  [\"591\u00d75=2955\", \"208\u00d78=1664\"],\n  [\"230\u00d76=1380\", \"140\u00d74=560\"],\n  [\"540\u00d72=1080\", \"179\u00d73=537\"],\n  [\"598\u00d76=3588\", \"812\u00d79=7308\"],\n  [\"805\u00d75=4025\", \"166\u00d72=332\"],\n  [\"226\u00d76=1356\", \"614\u00d73=1842\"],\n  [\"445\u00d75=2225\", \"330\u00d75=1650\"],\n  [\"478\u00d76=2868\", \"838\u00d78=6704\"],\n  [\"718\u00d78=5744\", \"374\u00d77=2618\"],\n  [\"499\u00d78=3992\", \"266\u00d79=2394\"],\n  [\"609\u00d75=3045\", \"804\u00d76=4824\"],\n  [\"228\u00d79=2052\", \"545\u00d77=3815\"],\n  [\"580\u00d75=2900\", \"397\u00d75=1985\"],\n  [\"460\u00d78=3680\", \"973\u00d73=2919\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell with its new value.\n# Every original text value in this document is unique, so an exact\n# match-case Find/Replace (wdReplaceAll) is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-07-08 Monday\", \"2024-07-09 Tuesday\"),\n    @(\"909\u00d77=6363\", \"614\u00d77=4298\"),\n    @(\"887\u00d74=3548\", \"721\u00d79=6489\"),\n    @(\"128\u00d77=896\", \"731\u00d72=1462\"),\n    @(\"617\u00d78=4936\", \"384\u00d77=2688\"),\n    @(\"559\u00d74=2236\", \"332\u00d73=996\"),\n    @(\"150\u00d72=300\", \"790\u00d78=6320\"),\n    @(\"584\u00d77=4088\", \"280\u00d74=1120\"),\n    @(\"995\u00d73=2985\", \"384\u00d72=768\"),\n    @(\"646\u00d78=5168\", \"370\u00d73=1110\"),\n    @(\"490\u00d73=1470\", \"761\u00d78=6088\"),\n    @(\"938\u00d72=1876\", \"743\u00d78=5944\"),\n    @(\"591\u00d75=2955\", \"208\u00d78=1664\"),\n    @(\"230\u00d76=1380\", \"140\u00d74=560\"),\n    @(\"540\u00d72=1080\", \"179\u00d73=537\"),\n    @(\"598\u00d76=3588\", \"812\u00d79=7308\"),\n    @(\"805\u00d75=4025\", \"166\u00d72=332\"),\n    @(\"226\u00d76=1356\", \"614\u00d73=1842\"),\n    @(\"445\u00d75=2225\", \"330\u00d75=1650\"),\n    @(\"478\u00d76=2868\", \"838\u00d78=6704\"),\n    @(\"718\u00d78=5744\", \"374\u00d77=2618\"),\n    @(\"499\u00d78=3992\", \"266\u00d79=2394\"),\n    @(\"609\u00d75=3045\", \"804\u00d76=4824\"),\n    @(\"228\u00d79=2052\", \"545\u00d77=3815\"),\n    @(\"580\u00d75=2900\", \"397\u00d75=1985\"),\n    @(\"460\u00d78=3680\", \"973\u00d73=2919\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
